# Update the "department" column (C) on the "courses" sheet.
# Previously every course row used the single generic department value
# "SHELDON SCHOOL OF HOSPITALITY"; it is now replaced with the specific
# department that each course belongs to.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

$ws.Range("C2").Value = "Cookery"                # SIT30816 - Cert III Commercial Cookery
$ws.Range("C3").Value = "Cookery"                # SIT40516 - Cert IV Commercial Cookery (fast track)
$ws.Range("C4").Value = "Patisserie and Baking"  # SIT31016 - Cert III Patisserie
$ws.Range("C5").Value = "Patisserie and Baking"  # SIT40716 - Cert IV Patisserie (fast track)
$ws.Range("C6").Value = "Patisserie and Baking"  # FBP40217 - Cert IV Baking
$ws.Range("C7").Value = "Hospitality"            # SIT50416 - Diploma of Hospitality Management (fast track)
$ws.Range("C8").Value = "Travel and Tourism"     # SIT50116 - Diploma of Travel and Tourism Management
$ws.Range("C9").Value = "Travel and Tourism"     # SIT60116 - Advanced Diploma of Travel and Tourism Management
